$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.854221333333333
$ws.Range("H2").Value = 5.562664
$ws.Range("I2").Value = 0.03416002559055492
$ws.Range("J2").Value = 0.03416002559055492
$ws.Range("M2").Value = 200.005264
$ws.Range("N2").Value = 600.0157919999999
$ws.Range("O2").Value = 0.9465949791503665
$ws.Range("P2").Value = 0.9465949791503667
$ws.Range("Q2").Value = 370.8540272877653
$ws.Range("R2").Value = 3337.686245589888
$ws.Range("S2").Value = 0.03233570871166732
$ws.Range("T2").Value = 0.03233570871166733
$ws.Range("G3").Value = 1.854221333333333
$ws.Range("H3").Value = 5.562664
$ws.Range("I3").Value = 0.03416002559055492
$ws.Range("J3").Value = 0.03416002559055492
$ws.Range("O3").Value = 0.006425713585924051
$ws.Range("P3").Value = 0.006425713585924052
$ws.Range("Q3").Value = 2.517446018651555
$ws.Range("R3").Value = 22.657014167864
$ws.Range("S3").Value = 0.000219502540532742
$ws.Range("T3").Value = 0.000219502540532742
$ws.Range("G4").Value = 1.854221333333333
$ws.Range("H4").Value = 5.562664
$ws.Range("I4").Value = 0.03416002559055492
$ws.Range("J4").Value = 0.03416002559055492
$ws.Range("M4").Value = 9.926218666666667
$ws.Range("O4").Value = 0.04697930726370939
$ws.Range("P4").Value = 0.0469793072637094
$ws.Range("Q4").Value = 18.40540641106489
$ws.Range("R4").Value = 165.648657699584
$ws.Range("S4").Value = 0.001604814338354855
$ws.Range("T4").Value = 0.001604814338354856
$ws.Range("I5").Value = 0.8311547934421808
$ws.Range("J5").Value = 0.8311547934421808
$ws.Range("M5").Value = 200.005264
$ws.Range("N5").Value = 600.0157919999999
$ws.Range("O5").Value = 0.9465949791503665
$ws.Range("P5").Value = 0.9465949791503667
$ws.Range("Q5").Value = 9023.327621065644
$ws.Range("R5").Value = 81209.94858959079
$ws.Range("S5").Value = 0.7867669543691284
$ws.Range("T5").Value = 0.7867669543691285
$ws.Range("I6").Value = 0.8311547934421808
$ws.Range("J6").Value = 0.8311547934421808
$ws.Range("O6").Value = 0.006425713585924051
$ws.Range("P6").Value = 0.006425713585924052
$ws.Range("S6").Value = 0.00534076264822732
$ws.Range("T6").Value = 0.005340762648227321
$ws.Range("I7").Value = 0.8311547934421808
$ws.Range("J7").Value = 0.8311547934421808
$ws.Range("M7").Value = 9.926218666666667
$ws.Range("O7").Value = 0.04697930726370939
$ws.Range("P7").Value = 0.0469793072637094
$ws.Range("Q7").Value = 447.8258285625459
$ws.Range("R7").Value = 4030.432457062913
$ws.Range("S7").Value = 0.03904707642482512
$ws.Range("T7").Value = 0.03904707642482513
$ws.Range("I8").Value = 0.1346851809672642
$ws.Range("J8").Value = 0.1346851809672642
$ws.Range("M8").Value = 200.005264
$ws.Range("N8").Value = 600.0157919999999
$ws.Range("O8").Value = 0.9465949791503665
$ws.Range("P8").Value = 0.9465949791503667
$ws.Range("Q8").Value = 1462.192750566965
$ws.Range("R8").Value = 13159.73475510269
$ws.Range("S8").Value = 0.1274923160695708
$ws.Range("T8").Value = 0.1274923160695708
$ws.Range("I9").Value = 0.1346851809672642
$ws.Range("J9").Value = 0.1346851809672642
$ws.Range("O9").Value = 0.006425713585924051
$ws.Range("P9").Value = 0.006425713585924052
$ws.Range("S9").Value = 0.0008654483971639892
$ws.Range("T9").Value = 0.0008654483971639893
$ws.Range("I10").Value = 0.1346851809672642
$ws.Range("J10").Value = 0.1346851809672642
$ws.Range("M10").Value = 9.926218666666667
$ws.Range("O10").Value = 0.04697930726370939
$ws.Range("P10").Value = 0.0469793072637094
$ws.Range("Q10").Value = 72.56831487666489
$ws.Range("R10").Value = 653.114833889984
$ws.Range("S10").Value = 0.006327416500529409
$ws.Range("T10").Value = 0.006327416500529411
